$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Sample_0199"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(3, 1).Value = "Sample_0172"
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(4, 1).Value = "Sample_0024"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(5, 1).Value = "Sample_0133"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(6, 1).Value = "Sample_0214"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(7, 1).Value = "Sample_0245"
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(8, 1).Value = "Sample_0123"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(9, 1).Value = "Sample_0106"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(10, 1).Value = "Sample_0026"
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(11, 1).Value = "Sample_0229"
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(12, 1).Value = "Sample_0216"
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(13, 1).Value = "Sample_0023"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(14, 1).Value = "Sample_0221"
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(15, 1).Value = "Sample_0111"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(16, 1).Value = "Sample_0174"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(17, 1).Value = "Sample_0173"
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(18, 1).Value = "Sample_0205"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(19, 1).Value = "Sample_0016"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(20, 1).Value = "Sample_0128"
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(21, 1).Value = "Sample_0247"
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(22, 1).Value = "Sample_0171"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(23, 1).Value = "Sample_0151"
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(24, 1).Value = "Sample_0198"
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(25, 1).Value = "Sample_0092"
$ws.Cells.Item(25, 2).Value = 1
$ws.Cells.Item(26, 1).Value = "Sample_0143"
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(27, 1).Value = "Sample_0120"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 1).Value = "Sample_0188"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(29, 1).Value = "Sample_0233"
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(30, 1).Value = "Sample_0045"
$ws.Cells.Item(30, 2).Value = 1
$ws.Cells.Item(31, 1).Value = "Sample_0122"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(32, 1).Value = "Sample_0140"
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(33, 1).Value = "Sample_0179"
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(34, 1).Value = "Sample_0243"
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(35, 1).Value = "Sample_0035"
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(36, 1).Value = "Sample_0197"
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(37, 1).Value = "Sample_0153"
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(38, 1).Value = "Sample_0107"
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(39, 1).Value = "Sample_0027"
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(40, 1).Value = "Sample_0182"
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(41, 1).Value = "Sample_0134"
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 1).Value = "Sample_0012"
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 1).Value = "Sample_0118"
$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(44, 1).Value = "Sample_0196"
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(45, 1).Value = "Sample_0155"
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(46, 1).Value = "Sample_0228"
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(47, 1).Value = "Sample_0183"
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(48, 1).Value = "Sample_0225"
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(49, 1).Value = "Sample_0019"
$ws.Cells.Item(49, 2).Value = 1
$ws.Cells.Item(50, 1).Value = "Sample_0062"
$ws.Cells.Item(50, 2).Value = 1
$ws.Cells.Item(51, 1).Value = "Sample_0222"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(52, 1).Value = "Sample_0089"
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(53, 1).Value = "Sample_0119"
$ws.Cells.Item(53, 2).Value = 1
$ws.Cells.Item(54, 1).Value = "Sample_0195"
$ws.Cells.Item(54, 2).Value = 1
$ws.Cells.Item(55, 1).Value = "Sample_0136"
$ws.Cells.Item(55, 2).Value = 0
$ws.Cells.Item(56, 1).Value = "Sample_0178"
$ws.Cells.Item(56, 2).Value = 0
$ws.Cells.Item(57, 1).Value = "Sample_0184"
$ws.Cells.Item(57, 2).Value = 0
$ws.Cells.Item(58, 1).Value = "Sample_0044"
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(59, 1).Value = "Sample_0230"
$ws.Cells.Item(59, 2).Value = 1
$ws.Cells.Item(60, 1).Value = "Sample_0230"
$ws.Cells.Item(60, 2).Value = 1
